# Edit: add two new weekly columns for "15.–21. 3. 2021" and "22.–28. 3. 2021"
# to both the "data" and "pocetR" sheets, and refresh the "aktualizace" date in the
# footnote cells from 23. 3. 2021 to 7. 4. 2021.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("data")
$wsPocet = $wb.Worksheets.Item("pocetR")

$newWeek1 = "15.–21. 3. 2021"
$newWeek2 = "22.–28. 3. 2021"

# ---------------------------------------------------------------------------
# Sheet "data": existing week columns run D..AW (col 4..49); append AX (50) and AY (51)
# ---------------------------------------------------------------------------

# Copy header style/format from the last existing week column (AW1) onto the two new header cells
$wsData.Range("AW1").Copy()
$wsData.Range("AX1:AY1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsData.Cells.Item(1, 50).Value = $newWeek1
$wsData.Cells.Item(1, 51).Value = $newWeek2

$dataRows = @"
2,0.61,0.61
3,0.09,0.09
4,0.16,0.16
5,0.14,0.14
6,0.65,0.63
7,0.07000000000000001,0.08
8,0.13,0.13
9,0.15,0.16
10,0.28,0.31
11,0.11,0.1
12,0.58,0.55
13,0.03,0.04
14,0.78,0.77
15,0.07000000000000001,0.06
16,0.03,0.03
17,0.12,0.14
18,0.64,0.65
19,0.13,0.14
20,0.1,0.1
21,0.13,0.11
22,0.39,0.36
23,0.14,0.19
24,0.4,0.37
25,0.07000000000000001,0.08
26,0.65,0.64
27,0.07000000000000001,0.08
28,0.14,0.14
29,0.14,0.14
30,0.67,0.66
31,0.09,0.08
32,0.1,0.11
33,0.14,0.15
34,0.5600000000000001,0.57
35,0.1,0.09
36,0.19,0.18
37,0.15,0.16
38,0.44,0.47
39,0.15,0.17
40,0.32,0.27
41,0.09,0.09
42,0.59,0.58
43,0.08,0.06
44,0.25,0.25
45,0.08,0.11
46,0.65,0.65
47,0.08,0.1
48,0.13,0.12
49,0.14,0.13
50,0.52,0.54
51,0.12,0.12
52,0.15,0.14
53,0.21,0.2
54,0.63,0.64
55,0.08,0.07000000000000001
56,0.13,0.12
57,0.16,0.17
58,0.59,0.57
59,0.1,0.12
60,0.21,0.2
61,0.1,0.11
62,0.62,0.62
63,0.09,0.09
64,0.16,0.15
65,0.13,0.14
66,0.72,0.71
67,0.04,0.03
68,0.04,0.03
69,0.2,0.23
70,0.64,0.64
71,0.1,0.11
72,0.15,0.15
73,0.11,0.1
74,0.41,0.42
75,0.14,0.16
76,0.38,0.35
77,0.07000000000000001,0.07000000000000001
"@

foreach ($line in ($dataRows -split "`n")) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $axVal = [double]$parts[1]
    $ayVal = [double]$parts[2]
    $wsData.Cells.Item($r, 50).Value = $axVal
    $wsData.Cells.Item($r, 51).Value = $ayVal
}

# Refresh the footnote date
$wsData.Cells.Item(78, 1).Value = "Život během pandemie, Home office, % respondentů celkově a ve skupinách, aktualizace 7. 4. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR": existing week columns run C..AV (col 3..48); append AW (49) and AX (50)
# ---------------------------------------------------------------------------

$wsPocet.Range("AV1").Copy()
$wsPocet.Range("AW1:AX1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsPocet.Cells.Item(1, 49).Value = $newWeek1
$wsPocet.Cells.Item(1, 50).Value = $newWeek2

$pocetRows = @"
2,1142,1142
3,316,316
4,95,95
5,304,304
6,168,168
7,98,98
8,558,558
9,269,269
10,147,147
11,168,168
12,267,267
13,650,650
14,225,225
15,432,432
16,371,371
17,339,339
18,410,410
19,456,456
20,276,276
"@

foreach ($line in ($pocetRows -split "`n")) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $r = [int]$parts[0]
    $awVal = [double]$parts[1]
    $axVal = [double]$parts[2]
    $wsPocet.Cells.Item($r, 49).Value = $awVal
    $wsPocet.Cells.Item($r, 50).Value = $axVal
}

# Trailing placeholder row (21) also gains the two new (blank) cells
$wsPocet.Cells.Item(21, 49).Value = ""
$wsPocet.Cells.Item(21, 50).Value = ""

# Refresh the footnote date
$wsPocet.Cells.Item(21, 1).Value = "Život během pandemie, Home office, velikost dotázaného souboru celkově a ve skupinách, aktualizace 7. 4. 2021"
